$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proposed Indicators")
$ws.Rows.Item(6).Delete() | Out-Null
$ws.Range("A5").Select() | Out-Null
